$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "types" table (H column): was id/name/??/color with a blank row,
#     fix to id/product_id/name/color (migration fix) ---
$ws.Range("H6").Value = "product_id"
$ws.Range("H7").Value = "name"

# H8 was blank (style s=2); give it the same boxed-row style as the rest
# of the table (e.g. H5/H7, style index 5) and fill in "color".
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("H8").Value = "color"

# --- "products" table (F column): qty now terminates the box (like the
#     other tables' last row, e.g. F18 "received"), and the old F12
#     "qty" row becomes an empty cap row with just a top border ---
$ws.Range("F18").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Value = "qty"

$ws.Range("F12").Clear()
$ws.Range("F12").Borders(8).LineStyle = "Single"

$excel.CutCopyMode = $false
